# Update the "Förändrad" date column (C) for rows 2-11 from 45183 to 45184.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value = 45184
    }
}
